{"js": "// Word Developer Log edit:\n//  1. Several paragraphs had their text split across many runs, interleaved\n//     with <w:proofErr/> spell/grammar-check markers (from Word's proofing\n//     pass). Re-writing each such paragraph's text collapses it back down to\n//     a single plain run and drops the stale proofErr markers.\n//  2. The trailing empty paragraph at the end of the log is replaced with\n//     three new dev-log entries.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraphs whose runs/proofErr markers need to be collapsed back into a\n// single plain run. Matched by their current (already proofErr-stripped)\n// text so the script does not depend on fragile positional indices.\nconst normalizations = [\n  \"Started Board.py.  __init__ creates the window and sets the coordinate system, and creates the board rectangle.  draw() draws the board.  wait_for_click() waits for a click.  test() creates a board and draws it, then calls wait_for_click().\",\n  \"Added ability to detect button clicks.  test() now waits for Done to be clicked.\",\n  \"Added move_tile(), added code to move tile 1 to upper left.\",\n  \"Added code to __init__ to set initial positions, and added code to track tile positions.  Added move_to_blank(), added test code to move 12 to the blank.\",\n  \"Added set_board(), called from __init__.  This allows us to provide a randomized initial position.\",\n  \"Removed draw(), added draw calls to __init__.\"\n];\n\nfor (const targetText of normalizations) {\n  const match = paragraphs.items.find((p) => p.text === targetText);\n  if (match) {\n    match.insertText(targetText, \"Replace\");\n  }\n}\nawait context.sync();\n\n// Replace the final (empty) paragraph with three new log entries.\nconst newEntries = [\n  \"Removed button creation from Board, replaced with add_button().  The game itself now creates the buttons.  This way there is less hard-coded game logic in the Board class.\",\n  \"The Game class now provides the initial configuration to the Board class.\",\n  \"Added copious documentation.\"\n];\n\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items/text\");\nawait context.sync();\n\nconst lastParagraph = refreshed.items[refreshed.items.length - 1];\n\nlet insertAfter = lastParagraph;\nif (lastParagraph.text === \"\") {\n  // Use the existing empty paragraph for the first new entry, then append\n  // the rest after it, so the trailing blank paragraph is replaced (not\n  // just followed) by the new text, matching the target document.\n  insertAfter.insertText(newEntries[0], \"Replace\");\n  for (let i = 1; i < newEntries.length; i++) {\n    insertAfter = insertAfter.insertParagraph(newEntries[i], \"After\");\n  }\n} else {\n  for (const text of newEntries) {\n    insertAfter = insertAfter.insertParagraph(text, \"After\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word Developer Log edit:\n#  1. Several paragraphs had their text split across many runs, interleaved\n#     with proofErr spell/grammar-check markers (from Word's proofing pass).\n#     Find & Replace'ing each paragraph's full text (old text -> identical\n#     new text) collapses it back down to a single plain run and drops the\n#     stale proofErr markers.\n#  2. The trailing empty paragraph at the end of the log is replaced with\n#     three new dev-log entries.\n\n$d = $word.ActiveDocument\n\n# Paragraphs whose runs/proofErr markers need to be collapsed back into a\n# single plain run.\n$normalizations = @(\n  \"Started Board.py.  __init__ creates the window and sets the coordinate system, and creates the board rectangle.  draw() draws the board.  wait_for_click() waits for a click.  test() creates a board and draws it, then calls wait_for_click().\",\n  \"Added ability to detect button clicks.  test() now waits for Done to be clicked.\",\n  \"Added move_tile(), added code to move tile 1 to upper left.\",\n  \"Added code to __init__ to set initial positions, and added code to track tile positions.  Added move_to_blank(), added test code to move 12 to the blank.\",\n  \"Added set_board(), called from __init__.  This allows us to provide a randomized initial position.\",\n  \"Removed draw(), added draw calls to __init__.\"\n)\n\nforeach ($text in $normalizations) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null\n}\n\n# Replace the final (empty) paragraph with three new log entries: reuse the\n# existing empty trailing paragraph for the first new entry, then add two\n# more paragraphs after it.\n$newEntries = @(\n  \"Removed button creation from Board, replaced with add_button().  The game itself now creates the buttons.  This way there is less hard-coded game logic in the Board class.\",\n  \"The Game class now provides the initial configuration to the Board class.\",\n  \"Added copious documentation.\"\n)\n\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$lastParagraph.Range.InsertAfter($newEntries[0])\n\nfor ($i = 1; $i -lt $newEntries.Count; $i++) {\n    $tail = $d.Paragraphs($d.Paragraphs.Count)\n    $tail.Range.InsertParagraphAfter()\n    $newTail = $d.Paragraphs($d.Paragraphs.Count)\n    $newTail.Range.InsertAfter($newEntries[$i])\n}\n"}
